# Add season record columns (Wins, Losses, Ties) to the BOS_2011 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - mirror the style of the existing header cells (e.g. AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-51: every player row gets the team's season record
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 90  # AD
    $ws.Cells.Item($r, 31).Value = 72  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
